$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for column C (y_0_forecast) and column E (y_1_forecast), rows 4-19
$values = @{
    4  = @{ C = -1.0538757111625;      E = -0.03775702182785556 }
    5  = @{ C = 1.570920254756558;     E = 0.485781149991249 }
    6  = @{ C = 0.878183952843048;     E = 0.6274179984581574 }
    7  = @{ C = 0.4413109953541605;    E = 0.6943226143418357 }
    8  = @{ C = 0.9083436352771646;    E = 0.8175300924317952 }
    9  = @{ C = 1.93704103170067;      E = 0.8819846436028733 }
    10 = @{ C = 1.58745197360306;      E = 0.9174072565958813 }
    11 = @{ C = 1.276746817047392;     E = 0.9840692158344266 }
    12 = @{ C = 1.116346046342809;     E = 1.076386909629967 }
    13 = @{ C = 1.437659299153027;     E = 1.280675138581411 }
    14 = @{ C = -1.741236505435428;    E = -0.2374583495023508 }
    15 = @{ C = -3.281345655000223;    E = -0.1411395765833157 }
    16 = @{ C = 4.407642498961351;     E = 1.101211681542535 }
    17 = @{ C = -1.129337961135779;    E = 0.5661352122229735 }
    18 = @{ C = 0.2034622429862187;    E = 0.6701250445558804 }
    19 = @{ C = 0.9632017291179906;    E = 0.8649794511895736 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
